$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1_Highlevel View")
$cell = $ws.Range("A1")
$cell.Borders.Item(1).LineStyle = -4142
$cell.Borders.Item(2).LineStyle = -4142
$cell.Borders.Item(3).LineStyle = -4142
$cell.Borders.Item(4).LineStyle = -4142
Write-Output "done"
